$d = $word.ActiveDocument

# Remove the trailing ", sorting dates" text (the comma that ended the
# previous run plus the whole " sorting dates" run) so the sentence ends
# at "...adjusted JSON of new API".
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(", sorting dates", $false, $false, $false, $false, $false, `
              $true, 1, $false, "", 2)
